$d = $word.ActiveDocument

# Update the title/date paragraph
$d.Content.Find.Execute("2024-12-22 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-23 Monday", 2) | Out-Null

# Update table cells (20 rows x 5 cols) by position to avoid any text-collision ambiguity
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "83-18="
$t.Cell(1, 2).Range.Text = "40-7="
$t.Cell(1, 3).Range.Text = "19+27="
$t.Cell(1, 4).Range.Text = "70-32="
$t.Cell(1, 5).Range.Text = "38+55="
$t.Cell(2, 1).Range.Text = "13-4="
$t.Cell(2, 2).Range.Text = "36+7="
$t.Cell(2, 3).Range.Text = "16+59="
$t.Cell(2, 4).Range.Text = "44+9="
$t.Cell(2, 5).Range.Text = "44-26="
$t.Cell(3, 1).Range.Text = "53-15="
$t.Cell(3, 2).Range.Text = "44+27="
$t.Cell(3, 3).Range.Text = "3+29="
$t.Cell(3, 4).Range.Text = "26+59="
$t.Cell(3, 5).Range.Text = "72-67="
$t.Cell(4, 1).Range.Text = "91-56="
$t.Cell(4, 2).Range.Text = "75-38="
$t.Cell(4, 3).Range.Text = "19+9="
$t.Cell(4, 4).Range.Text = "19+12="
$t.Cell(4, 5).Range.Text = "82-77="
$t.Cell(5, 1).Range.Text = "5+6="
$t.Cell(5, 2).Range.Text = "80-69="
$t.Cell(5, 3).Range.Text = "21-5="
$t.Cell(5, 4).Range.Text = "24+37="
$t.Cell(5, 5).Range.Text = "44-37="
$t.Cell(6, 1).Range.Text = "42-14="
$t.Cell(6, 2).Range.Text = "35-19="
$t.Cell(6, 3).Range.Text = "57-8="
$t.Cell(6, 4).Range.Text = "57+39="
$t.Cell(6, 5).Range.Text = "47-9="
$t.Cell(7, 1).Range.Text = "47+14="
$t.Cell(7, 2).Range.Text = "76-27="
$t.Cell(7, 3).Range.Text = "79+2="
$t.Cell(7, 4).Range.Text = "69+8="
$t.Cell(7, 5).Range.Text = "7+38="
$t.Cell(8, 1).Range.Text = "70-42="
$t.Cell(8, 2).Range.Text = "82-47="
$t.Cell(8, 3).Range.Text = "28+18="
$t.Cell(8, 4).Range.Text = "35+49="
$t.Cell(8, 5).Range.Text = "8+15="
$t.Cell(9, 1).Range.Text = "2+89="
$t.Cell(9, 2).Range.Text = "51-26="
$t.Cell(9, 3).Range.Text = "81-17="
$t.Cell(9, 4).Range.Text = "3+18="
$t.Cell(9, 5).Range.Text = "61-32="
$t.Cell(10, 1).Range.Text = "73-26="
$t.Cell(10, 2).Range.Text = "30-8="
$t.Cell(10, 3).Range.Text = "85-26="
$t.Cell(10, 4).Range.Text = "48+35="
$t.Cell(10, 5).Range.Text = "18+24="
$t.Cell(11, 1).Range.Text = "57+39="
$t.Cell(11, 2).Range.Text = "93-34="
$t.Cell(11, 3).Range.Text = "3+69="
$t.Cell(11, 5).Range.Text = "5+16="
$t.Cell(12, 1).Range.Text = "47+5="
$t.Cell(12, 2).Range.Text = "93-38="
$t.Cell(12, 3).Range.Text = "73-57="
$t.Cell(12, 4).Range.Text = "61-55="
$t.Cell(12, 5).Range.Text = "17+49="
$t.Cell(13, 1).Range.Text = "9+38="
$t.Cell(13, 2).Range.Text = "6+27="
$t.Cell(13, 3).Range.Text = "52-17="
$t.Cell(13, 4).Range.Text = "9+23="
$t.Cell(13, 5).Range.Text = "15-6="
$t.Cell(14, 1).Range.Text = "61-38="
$t.Cell(14, 2).Range.Text = "7+4="
$t.Cell(14, 3).Range.Text = "46+35="
$t.Cell(14, 4).Range.Text = "29+27="
$t.Cell(14, 5).Range.Text = "68+19="
$t.Cell(15, 1).Range.Text = "71-7="
$t.Cell(15, 2).Range.Text = "78-9="
$t.Cell(15, 3).Range.Text = "38+6="
$t.Cell(15, 4).Range.Text = "29+67="
$t.Cell(15, 5).Range.Text = "77-59="
$t.Cell(16, 1).Range.Text = "53-15="
$t.Cell(16, 2).Range.Text = "47+48="
$t.Cell(16, 3).Range.Text = "16+59="
$t.Cell(16, 4).Range.Text = "74-49="
$t.Cell(16, 5).Range.Text = "95-79="
$t.Cell(17, 1).Range.Text = "12+19="
$t.Cell(17, 2).Range.Text = "23-7="
$t.Cell(17, 3).Range.Text = "27+37="
$t.Cell(17, 4).Range.Text = "61-48="
$t.Cell(17, 5).Range.Text = "46+26="
$t.Cell(18, 1).Range.Text = "88+4="
$t.Cell(18, 2).Range.Text = "10-9="
$t.Cell(18, 3).Range.Text = "84-29="
$t.Cell(18, 4).Range.Text = "46+7="
$t.Cell(18, 5).Range.Text = "6+26="
$t.Cell(19, 1).Range.Text = "73+8="
$t.Cell(19, 2).Range.Text = "32-19="
$t.Cell(19, 3).Range.Text = "53-37="
$t.Cell(19, 4).Range.Text = "57-38="
$t.Cell(19, 5).Range.Text = "17+74="
$t.Cell(20, 1).Range.Text = "85-6="
$t.Cell(20, 2).Range.Text = "6+55="
$t.Cell(20, 3).Range.Text = "19+35="
$t.Cell(20, 4).Range.Text = "62-48="
$t.Cell(20, 5).Range.Text = "73-24="
